# Fruta / hortaliza, semanal
# Insert this week's new price observation as a new record in the
# "Papa" (potato) price table. Excel's row-insert naturally pushes the
# existing record that was at row 520 (and everything below it) down by
# one row, which is exactly the shift seen between before/after states.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("520:520").Insert()

$ws.Range("A520").Value = 8
$ws.Range("B520").Value = "Terminal La Palmera de La Serena"
$ws.Range("C520").Value = "Coquimbo"
$ws.Range("D520").Value = 44931
$ws.Range("E520").Value = 4
$ws.Range("F520").Value = 100114001
$ws.Range("G520").Value = "Papa"
$ws.Range("H520").Value = "Cardinal"
$ws.Range("I520").Value = "1a (cosecha)"
$ws.Range("J520").Value = 2200
$ws.Range("K520").Value = 11500
$ws.Range("L520").Value = 12000
$ws.Range("M520").Value = 11750
$ws.Range("N520").Value = '$/saco 25 kilos'
$ws.Range("O520").Value = "Provincia del Elquí"
$ws.Range("P520").Value = 470
$ws.Range("Q520").Value = 25
$ws.Range("R520").Value = "Hortaliza"
